$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

$ws.Range("D2").Value = '26.963.74'
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").Value = '1.560.59'
$ws.Range("E3").Value = '  +0.42%  '
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").Value = '207.13'
$ws.Range("E5").Value = '  +0.08%  '
$ws.Range("E6").Value = '  +0.46%  '
$ws.Range("E7").Value = '  -0.23%  '
$ws.Range("D8").Value = '22.12'
$ws.Range("E8").Value = '  +1.65%  '
$ws.Range("E9").Value = '  +0.18%  '
$ws.Range("E10").Value = '  +2.09%  '
$ws.Range("E11").Value = '  -0.08%  '
$ws.Range("D12").Value = '1.783.45'
$ws.Range("D13").Value = '1.566.31'
$ws.Range("E13").Value = '  +0.73%  '
$ws.Range("E14").Value = '  +0.58%  '
$ws.Range("E15").Value = '  +0.47%  '
$ws.Range("D16").Value = '62.16'
$ws.Range("E16").Value = '  +0.77%  '
$ws.Range("D17").Value = '26.974.01'
$ws.Range("E17").Value = '  +0.16%  '
$ws.Range("D18").Value = '217.13'
$ws.Range("E18").Value = '  +0.00%  '
$ws.Range("D19").Value = '0.0₃0702'
$ws.Range("E19").Value = '  +2.09%  '
$ws.Range("D20").Value = '7.35'
$ws.Range("E20").Value = '  +1.76%  '
$ws.Range("E21").Value = '  -0.19%  '
$ws.Range("D22").Value = '4.11'
$ws.Range("E22").Value = '  +1.69%  '
$ws.Range("D23").Value = '9.19'
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("E24").Value = '  -1.29%  '
$ws.Range("D25").Value = '153.29'
$ws.Range("E25").Value = '  -0.35%  '
$ws.Range("E26").Value = '  +0.75%  '
$ws.Range("D27").Value = '15.08'
$ws.Range("E27").Value = '  +1.42%  '
$ws.Range("E29").Value = '  -0.30%  '
$ws.Range("D30").Value = '0.0470'
$ws.Range("E30").Value = '  +0.59%  '
$ws.Range("E31").Value = '  +1.37%  '
$ws.Range("E32").Value = '  +0.68%  '
$ws.Range("D33").Value = '3.12'
$ws.Range("E33").Value = '  +3.18%  '
$ws.Range("D34").Value = '1.421.41'
$ws.Range("E34").Value = '  -0.20%  '
$ws.Range("E35").Value = '  +3.03%  '
$ws.Range("E36").Value = '  +8.56%  '
$ws.Range("E37").Value = '  +1.53%  '
$ws.Range("E38").Value = '  +0.71%  '
$ws.Range("E39").Value = '  +3.03%  '
$ws.Range("E40").Value = '  -0.05%  '
$ws.Range("E41").Value = '  -0.24%  '
$ws.Range("E42").Value = '  +0.25%  '
$ws.Range("D43").Value = '2.33'
$ws.Range("E43").Value = '  +2.68%  '
$ws.Range("E44").Value = '  +2.08%  '
$ws.Range("D45").Value = '64.84'
$ws.Range("E45").Value = '  +1.68%  '
$ws.Range("D46").Value = '1.74'
$ws.Range("E46").Value = '  +0.67%  '
$ws.Range("D47").Value = '1.696.91'
$ws.Range("E47").Value = '  +0.43%  '
$ws.Range("D48").Value = '87.45'
$ws.Range("E48").Value = '  +1.45%  '
$ws.Range("E49").Value = '  -0.29%  '
$ws.Range("E50").Value = '  +0.04%  '
$ws.Range("E51").Value = '  -0.22%  '
